$d = $word.ActiveDocument
$nbsp = [char]160

# --- Change 1: reset the "smarthosting" bookmark id from 1 to 0 ---
# Word COM does not allow setting bookmark IDs directly; deleting and
# re-adding the bookmark causes Word to renumber bookmark ids starting at 0.
$bm = $d.Bookmarks("smarthosting")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("smarthosting", $bmRange)

# --- Change 2: drop the CryptoBridge / HitBTC exchange links ---
# Remove the two HYPERLINK fields whose visible text is "CryptoBridge" and
# "HitBTC" (this deletes the field begin/instrText/separate/end runs too).
$fieldsToDelete = @()
foreach ($f in $d.Fields) {
    if ($f.Result.Text -eq "CryptoBridge" -or $f.Result.Text -eq "HitBTC") {
        $fieldsToDelete += $f
    }
}
foreach ($f in $fieldsToDelete) {
    $f.Delete()
}

# Remove the now-dangling " such as<nbsp>" text that used to introduce the
# two links, leaving "...obtained from exchanges".
$r1 = $d.Content
[void]$r1.Find.Execute(" such as" + $nbsp, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Text = ""

# Remove the leftover "," + nbsp that used to separate the two links
# (leaving just the trailing "." that was already after "HitBTC").
$r2 = $d.Content
[void]$r2.Find.Execute("," + $nbsp, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Text = ""
